# Apply the 2020-07-21 data refresh to the "Fonds de solidarite volet 1"
# sheet: update columns C (nombre_aides) and D (montant_total) for the
# rows listed below, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @"
2|35199|50977093
3|86070|126325140
4|29509|43739115
5|8135|12097018
6|1752|2605096
7|131|191593
11|38674|52612658
12|9078|13137536
13|24790|36383432
14|7900|11734756
15|2011|2991976
16|372|547123
17|29|43500
18|6|9000
19|9584|12735739
20|12657|18290338
21|30278|44486054
22|9809|14588251
23|2482|3693263
24|443|658845
25|29|42953
26|11047|14815702
27|7209|10449697
28|21416|31456789
29|7452|11092709
30|1839|2746986
31|301|448915
32|26|38893
33|7844|10398358
34|2831|4083581
35|7095|10370149
36|2873|4253788
37|767|1143263
40|2190|2952114
41|16221|23476842
42|48573|71264361
43|18156|26975421
44|5270|7852784
45|1054|1572219
46|54|80191
49|15681|20939342
50|1731|2513982
51|6137|9033194
52|2120|3166750
53|698|1042305
54|161|239111
56|5628|7772522
57|714|1046540
58|1819|2696626
59|733|1091953
60|253|379258
61|54|81000
62|14|21000
63|1058|1501861
64|14442|20877295
65|42678|62497519
66|15043|22368115
67|4343|6468716
68|842|1252596
69|68|100689
71|14311|18943167
72|46677|67968779
73|134630|198482495
74|58945|87865801
75|18711|27962419
76|4203|6281120
77|226|334170
78|18|25905
81|4|6000
83|46135|63103026
84|4170|6047631
85|10740|15786309
86|3667|5465993
87|1274|1904913
88|264|393512
91|4838|6523494
92|1426|2062583
93|4684|6901416
94|1803|2688403
95|641|960641
96|159|237613
98|6|9000
99|3106|4121729
100|529|788964
101|276|412165
102|93|139500
103|38|57000
104|18|27000
105|10168|14773994
106|27984|41136250
107|9363|13925548
108|2549|3800910
109|437|652982
112|9247|12258927
113|28576|41249120
114|63029|92320248
115|20421|30363779
116|5721|8527458
117|1023|1530493
118|60|87420
121|24352|32605179
122|33305|48132950
123|72348|105902791
124|22589|33540084
125|5976|8885761
126|1094|1627646
130|29463|39244320
131|12475|18070229
132|30903|45422752
133|11010|16361344
134|2771|4132791
135|445|661490
138|10226|13695035
139|32557|47058800
140|76638|112359980
141|23079|34314224
142|5963|8900822
143|1296|1930686
144|69|103130
146|27344|37047653
"@

$lines = $updates -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $nombreAides = [int]$parts[1]
    $montantTotal = [int]$parts[2]

    $ws.Cells.Item($row, 3).Value = $nombreAides
    $ws.Cells.Item($row, 4).Value = $montantTotal
}
